$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = "'63.181.57"
$ws.Cells.Item(2, 4).Style = "Normal"
$ws.Cells.Item(2, 5).Value = "'  +0.62%  "
$ws.Cells.Item(2, 5).Style = "Normal"
$ws.Cells.Item(3, 4).Value = "'2.569.19"
$ws.Cells.Item(3, 4).Style = "Normal"
$ws.Cells.Item(3, 5).Value = "'  +1.04%  "
$ws.Cells.Item(3, 5).Style = "Normal"
$ws.Cells.Item(4, 5).Value = "'  +0.01%  "
$ws.Cells.Item(4, 5).Style = "Normal"
$ws.Cells.Item(5, 4).Value = "'585.25"
$ws.Cells.Item(5, 4).Style = "Normal"
$ws.Cells.Item(5, 5).Value = "'  +3.12%  "
$ws.Cells.Item(5, 5).Style = "Normal"
$ws.Cells.Item(6, 4).Value = "'147.52"
$ws.Cells.Item(6, 4).Style = "Normal"
$ws.Cells.Item(6, 5).Value = "'  +0.43%  "
$ws.Cells.Item(6, 5).Style = "Normal"
$ws.Cells.Item(7, 5).Value = "'  +0.00%  "
$ws.Cells.Item(7, 5).Style = "Normal"
$ws.Cells.Item(8, 4).Value = "'0.603"
$ws.Cells.Item(8, 4).Style = "Normal"
$ws.Cells.Item(8, 5).Value = "'  +3.20%  "
$ws.Cells.Item(8, 5).Style = "Normal"
$ws.Cells.Item(9, 5).Value = "'  +3.76%  "
$ws.Cells.Item(9, 5).Style = "Normal"
$ws.Cells.Item(10, 4).Value = "'5.64"
$ws.Cells.Item(10, 4).Style = "Normal"
$ws.Cells.Item(10, 5).Value = "'  +0.72%  "
$ws.Cells.Item(10, 5).Style = "Normal"
$ws.Cells.Item(11, 5).Value = "'  +0.08%  "
$ws.Cells.Item(11, 5).Style = "Normal"
$ws.Cells.Item(12, 4).Value = "'0.357"
$ws.Cells.Item(12, 4).Style = "Normal"
$ws.Cells.Item(12, 5).Value = "'  +1.51%  "
$ws.Cells.Item(12, 5).Style = "Normal"
$ws.Cells.Item(13, 4).Value = "'27.50"
$ws.Cells.Item(13, 4).Style = "Normal"
$ws.Cells.Item(13, 5).Value = "'  +1.09%  "
$ws.Cells.Item(13, 5).Style = "Normal"
$ws.Cells.Item(14, 4).Value = "'3.029.53"
$ws.Cells.Item(14, 4).Style = "Normal"
$ws.Cells.Item(14, 5).Value = "'  +1.09%  "
$ws.Cells.Item(14, 5).Style = "Normal"
$ws.Cells.Item(15, 4).Value = "'63.132.18"
$ws.Cells.Item(15, 4).Style = "Normal"
$ws.Cells.Item(15, 5).Value = "'  +0.53%  "
$ws.Cells.Item(15, 5).Style = "Normal"
$ws.Cells.Item(16, 5).Value = "'  +3.90%  "
$ws.Cells.Item(16, 5).Style = "Normal"
$ws.Cells.Item(17, 4).Value = "'2.597.44"
$ws.Cells.Item(17, 4).Style = "Normal"
$ws.Cells.Item(17, 5).Value = "'  +2.15%  "
$ws.Cells.Item(17, 5).Style = "Normal"
$ws.Cells.Item(18, 4).Value = "'11.37"
$ws.Cells.Item(18, 4).Style = "Normal"
$ws.Cells.Item(18, 5).Value = "'  -0.59%  "
$ws.Cells.Item(18, 5).Style = "Normal"
$ws.Cells.Item(19, 4).Value = "'343.94"
$ws.Cells.Item(19, 4).Style = "Normal"
$ws.Cells.Item(19, 5).Value = "'  +2.21%  "
$ws.Cells.Item(19, 5).Style = "Normal"
$ws.Cells.Item(20, 4).Value = "'4.43"
$ws.Cells.Item(20, 4).Style = "Normal"
$ws.Cells.Item(20, 5).Value = "'  +3.36%  "
$ws.Cells.Item(20, 5).Style = "Normal"
$ws.Cells.Item(21, 5).Value = "'  +2.35%  "
$ws.Cells.Item(21, 5).Style = "Normal"
$ws.Cells.Item(22, 5).Value = "'  +0.05%  "
$ws.Cells.Item(22, 5).Style = "Normal"
$ws.Cells.Item(23, 5).Value = "'  -3.77%  "
$ws.Cells.Item(23, 5).Style = "Normal"
$ws.Cells.Item(24, 4).Value = "'66.87"
$ws.Cells.Item(24, 4).Style = "Normal"
$ws.Cells.Item(24, 5).Value = "'  +2.49%  "
$ws.Cells.Item(24, 5).Style = "Normal"
$ws.Cells.Item(25, 4).Value = "'2.697.14"
$ws.Cells.Item(25, 4).Style = "Normal"
$ws.Cells.Item(25, 5).Value = "'  +1.04%  "
$ws.Cells.Item(25, 5).Style = "Normal"
$ws.Cells.Item(26, 5).Value = "'  +1.30%  "
$ws.Cells.Item(26, 5).Style = "Normal"
$ws.Cells.Item(27, 5).Value = "'  +0.86%  "
$ws.Cells.Item(27, 5).Style = "Normal"
$ws.Cells.Item(28, 4).Value = "'8.16"
$ws.Cells.Item(28, 4).Style = "Normal"
$ws.Cells.Item(28, 5).Value = "'  +12.55%  "
$ws.Cells.Item(28, 5).Style = "Normal"
$ws.Cells.Item(29, 4).Value = "'8.51"
$ws.Cells.Item(29, 4).Style = "Normal"
$ws.Cells.Item(29, 5).Value = "'  +1.64%  "
$ws.Cells.Item(29, 5).Style = "Normal"
$ws.Cells.Item(30, 2).Value = "'Binance-PegBSC-USD"
$ws.Cells.Item(30, 2).Style = "Normal"
$ws.Cells.Item(30, 3).Value = "'https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd"
$ws.Cells.Item(30, 3).Style = "Normal"
$ws.Cells.Item(30, 4).Value = "'1.00"
$ws.Cells.Item(30, 4).Style = "Normal"
$ws.Cells.Item(30, 5).Value = "'  -0.02%  "
$ws.Cells.Item(30, 5).Style = "Normal"
$ws.Cells.Item(31, 2).Value = "'SuiNetwork"
$ws.Cells.Item(31, 2).Style = "Normal"
$ws.Cells.Item(31, 3).Value = "'https://coinranking.com/coin/3xJluUMvp+suinetwork-sui"
$ws.Cells.Item(31, 3).Style = "Normal"
$ws.Cells.Item(31, 4).Value = "'1.48"
$ws.Cells.Item(31, 4).Style = "Normal"
$ws.Cells.Item(31, 5).Value = "'  -1.17%  "
$ws.Cells.Item(31, 5).Style = "Normal"
$ws.Cells.Item(32, 5).Value = "'  +7.99%  "
$ws.Cells.Item(32, 5).Style = "Normal"
$ws.Cells.Item(33, 4).Value = "'0.0₃0828"
$ws.Cells.Item(33, 4).Style = "Normal"
$ws.Cells.Item(33, 5).Value = "'  +2.47%  "
$ws.Cells.Item(33, 5).Style = "Normal"
$ws.Cells.Item(34, 4).Value = "'466.22"
$ws.Cells.Item(34, 4).Style = "Normal"
$ws.Cells.Item(34, 5).Value = "'  +13.55%  "
$ws.Cells.Item(34, 5).Style = "Normal"
$ws.Cells.Item(35, 4).Value = "'1.63"
$ws.Cells.Item(35, 4).Style = "Normal"
$ws.Cells.Item(35, 5).Value = "'  +3.78%  "
$ws.Cells.Item(35, 5).Style = "Normal"
$ws.Cells.Item(36, 4).Value = "'176.08"
$ws.Cells.Item(36, 4).Style = "Normal"
$ws.Cells.Item(36, 5).Value = "'  -0.84%  "
$ws.Cells.Item(36, 5).Style = "Normal"
$ws.Cells.Item(37, 5).Value = "'  +2.50%  "
$ws.Cells.Item(37, 5).Style = "Normal"
$ws.Cells.Item(38, 4).Value = "'19.24"
$ws.Cells.Item(38, 4).Style = "Normal"
$ws.Cells.Item(38, 5).Value = "'  +1.59%  "
$ws.Cells.Item(38, 5).Style = "Normal"
$ws.Cells.Item(39, 5).Value = "'  +4.53%  "
$ws.Cells.Item(39, 5).Style = "Normal"
$ws.Cells.Item(40, 5).Value = "'  +0.02%  "
$ws.Cells.Item(40, 5).Style = "Normal"
$ws.Cells.Item(41, 5).Value = "'  -0.29%  "
$ws.Cells.Item(41, 5).Style = "Normal"
$ws.Cells.Item(42, 5).Value = "'  +0.03%  "
$ws.Cells.Item(42, 5).Style = "Normal"
$ws.Cells.Item(43, 5).Value = "'  -0.76%  "
$ws.Cells.Item(43, 5).Style = "Normal"
$ws.Cells.Item(44, 4).Value = "'3.82"
$ws.Cells.Item(44, 4).Style = "Normal"
$ws.Cells.Item(44, 5).Value = "'  +2.27%  "
$ws.Cells.Item(44, 5).Style = "Normal"
$ws.Cells.Item(45, 5).Value = "'  +0.94%  "
$ws.Cells.Item(45, 5).Style = "Normal"
$ws.Cells.Item(46, 4).Value = "'0.0549"
$ws.Cells.Item(46, 4).Style = "Normal"
$ws.Cells.Item(46, 5).Value = "'  +5.93%  "
$ws.Cells.Item(46, 5).Style = "Normal"
$ws.Cells.Item(47, 4).Value = "'0.611"
$ws.Cells.Item(47, 4).Style = "Normal"
$ws.Cells.Item(47, 5).Value = "'  +1.12%  "
$ws.Cells.Item(47, 5).Style = "Normal"
$ws.Cells.Item(48, 4).Value = "'0.0980"
$ws.Cells.Item(48, 4).Style = "Normal"
$ws.Cells.Item(48, 5).Value = "'  +2.37%  "
$ws.Cells.Item(48, 5).Style = "Normal"
$ws.Cells.Item(49, 4).Value = "'0.0240"
$ws.Cells.Item(49, 4).Style = "Normal"
$ws.Cells.Item(49, 5).Value = "'  +1.90%  "
$ws.Cells.Item(49, 5).Style = "Normal"
$ws.Cells.Item(50, 4).Value = "'1.74"
$ws.Cells.Item(50, 4).Style = "Normal"
$ws.Cells.Item(50, 5).Value = "'  -1.48%  "
$ws.Cells.Item(50, 5).Style = "Normal"
$ws.Cells.Item(51, 5).Value = "'  -0.11%  "
$ws.Cells.Item(51, 5).Style = "Normal"
